# Update row 2 (300913.SZ cash-flow record) to the prior-year (2019-09-30) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-09-30 00:00:00"
$ws.Range("O2").Value = 68041416.11
$ws.Range("P2").Value = 377.9480404127
$ws.Range("Q2").Value = 819157013.5599999
$ws.Range("R2").Value = 4550.1520362947
$ws.Range("S2").Value = 65967001.56
$ws.Range("T2").Value = 366.4253391081
$ws.Range("U2").Value = -27791815.4
$ws.Range("V2").Value = -154.3745378985
$ws.Range("X2").Value = 0.4443741027
$ws.Range("Y2").Value = 27881815.4
$ws.Range("Z2").Value = 154.8744587641
$ws.Range("AA2").Value = -59153365.7
$ws.Range("AB2").Value = -328.5777975871
$ws.Range("AC2").Value = -18002849.29
$ws.Range("AD2").ClearContents()
